$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 414. This shifts existing rows 414-498 down to 415-499,
# matching the diff (old row 414 contents effectively move to row 415, ..., old row 498
# moves to new row 499), and the sheet dimension grows from T498 to T499.
$ws.Rows.Item(414).Insert()

# Populate the newly inserted row 414 with the new weekly data point.
$ws.Cells.Item(414, 1).Value = 4
$ws.Cells.Item(414, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(414, 3).Value = "Los Lagos"
$ws.Cells.Item(414, 4).Value = 44711
$ws.Cells.Item(414, 5).Value = 10
$ws.Cells.Item(414, 6).Value = "Fruta"
$ws.Cells.Item(414, 7).Value = 100102
$ws.Cells.Item(414, 8).Value = "Cítricos"
$ws.Cells.Item(414, 9).Value = 100102003
$ws.Cells.Item(414, 10).Value = "Limón"
$ws.Cells.Item(414, 11).Value = "Sin especificar"
$ws.Cells.Item(414, 12).Value = "1a plateado"
$ws.Cells.Item(414, 13).Value = 750
$ws.Cells.Item(414, 14).Value = 12000
$ws.Cells.Item(414, 15).Value = 14000
$ws.Cells.Item(414, 16).Value = 13333
$ws.Cells.Item(414, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(414, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(414, 19).Value = 741
$ws.Cells.Item(414, 20).Value = 18
